$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 8.684023833688684
$ws.Range("C2").Value = 0.9652823463268305
$ws.Range("D2").Value = 0.2857124914840554
$ws.Range("E2").Value = 0.03637896011827024
$ws.Range("G2").Value = 0.002840580847411919
$ws.Range("J2").Value = 0.02730041464540989
$ws.Range("L2").Value = 0.7165019072303807
$ws.Range("M2").Value = 1.61067766300998
$ws.Range("N2").Value = 6.210132110822883
$ws.Range("B3").Value = 8.568314884328572
$ws.Range("C3").Value = 0.9328352151426316
$ws.Range("D3").Value = 0.2865386693235905
$ws.Range("E3").Value = 0.03594057371530113
$ws.Range("G3").Value = 0.002849739966127345
$ws.Range("J3").Value = 0.02399673562509719
$ws.Range("L3").Value = 0.7150013713199428
$ws.Range("M3").Value = 1.594528402402908
$ws.Range("N3").Value = 6.135725961313085
$ws.Range("B4").Value = 8.502350826202075
$ws.Range("C4").Value = 0.9135852595194365
$ws.Range("D4").Value = 0.2871470892178607
$ws.Range("E4").Value = 0.03566638692483348
$ws.Range("G4").Value = 0.002855650285851822
$ws.Range("J4").Value = 0.02197906809903571
$ws.Range("L4").Value = 0.7144141970595967
$ws.Range("M4").Value = 1.585511455756546
$ws.Range("N4").Value = 6.090505686532396
$ws.Range("B5").Value = 8.476740659082566
$ws.Range("C5").Value = 0.905908033960543
$ws.Range("D5").Value = 0.2874204609057358
$ws.Range("E5").Value = 0.03555337014184623
$ws.Range("G5").Value = 0.002858131153120567
$ws.Range("J5").Value = 0.02115939871460171
$ws.Range("L5").Value = 0.7142587475192244
$ws.Range("M5").Value = 1.582061910512728
$ws.Range("N5").Value = 6.072192350085999
$ws.Range("B6").Value = 8.47256465463829
$ws.Range("C6").Value = 0.9046432862949132
$ws.Range("D6").Value = 0.287467390419387
$ws.Range("E6").Value = 0.03553452563041315
$ws.Range("G6").Value = 0.002858547478595263
$ws.Range("J6").Value = 0.02102344210828022
$ws.Range("L6").Value = 0.7142379926539348
$ws.Range("M6").Value = 1.58150267098916
$ws.Range("N6").Value = 6.069158258115124
$ws.Range("B7").Value = 8.502000302654153
$ws.Range("C7").Value = 0.9134810468719365
$ws.Range("D7").Value = 0.2871506730113111
$ws.Range("E7").Value = 0.03566486796675505
$ws.Range("G7").Value = 0.002855683450336078
$ws.Range("J7").Value = 0.02196800366395735
$ws.Range("L7").Value = 0.714411761446911
$ws.Range("M7").Value = 1.585464024617323
$ws.Range("N7").Value = 6.090258247301051
$ws.Range("B8").Value = 8.643068286081075
$ws.Range("C8").Value = 0.9539538040020261
$ws.Range("D8").Value = 0.2859763658812255
$ws.Range("E8").Value = 0.03622883199783899
$ws.Range("G8").Value = 0.002843679612476921
$ws.Range("J8").Value = 0.02615895934365398
$ws.Range("L8").Value = 0.7159150431171781
$ws.Range("M8").Value = 1.604922198631236
$ws.Range("N8").Value = 6.184378598641558
$ws.Range("B9").Value = 8.960377024296463
$ws.Range("C9").Value = 1.038749495267325
$ws.Range("D9").Value = 0.2844762193084307
$ws.Range("E9").Value = 0.03729600832129609
$ws.Range("G9").Value = 0.00282240030337142
$ws.Range("J9").Value = 0.03447130227318951
$ws.Range("L9").Value = 0.7215247144158354
$ws.Range("M9").Value = 1.650264339412345
$ws.Range("N9").Value = 6.372774777782098
$ws.Range("B10").Value = 9.218838387146775
$ws.Range("C10").Value = 1.104494173377134
$ws.Range("D10").Value = 0.2838639403723846
$ws.Range("E10").Value = 0.03805794092482628
$ws.Range("G10").Value = 0.002808124964157705
$ws.Range("J10").Value = 0.04064761149138008
$ws.Range("L10").Value = 0.7272851025985148
$ws.Range("M10").Value = 1.688038890224306
$ws.Range("N10").Value = 6.513722246923294
$ws.Range("B11").Value = 9.342043240716407
$ws.Range("C11").Value = 1.135181895815435
$ws.Range("D11").Value = 0.2836919495005858
$ws.Range("E11").Value = 0.03840009330863303
$ws.Range("G11").Value = 0.002801921587992334
$ws.Range("J11").Value = 0.04347508796737998
$ws.Range("L11").Value = 0.7302653539321255
$ws.Range("M11").Value = 1.70621104569652
$ws.Range("N11").Value = 6.578437523767093
$ws.Range("B12").Value = 9.389516900625949
$ws.Range("C12").Value = 1.146917143254541
$ws.Range("D12").Value = 0.2836421554894883
$ws.Range("E12").Value = 0.03852904345266772
$ws.Range("G12").Value = 0.002799613987961212
$ws.Range("J12").Value = 0.04454855726537232
$ws.Range("L12").Value = 0.7314459416278396
$ws.Range("M12").Value = 1.713235962676208
$ws.Range("N12").Value = 6.603032846064821
$ws.Range("B13").Value = 9.379256039885263
$ws.Range("C13").Value = 1.144384619420748
$ws.Range("D13").Value = 0.2836521972167816
$ws.Range("E13").Value = 0.03850129870500751
$ws.Range("G13").Value = 0.002800109130550492
$ws.Range("J13").Value = 0.04431724023270078
$ws.Range("L13").Value = 0.7311893621030379
$ws.Range("M13").Value = 1.711716615454961
$ws.Range("N13").Value = 6.597731792407217
$ws.Range("B14").Value = 9.345932466809245
$ws.Range("C14").Value = 1.13614505413841
$ws.Range("D14").Value = 0.2836875454826071
$ws.Range("E14").Value = 0.03841071430582055
$ws.Range("G14").Value = 0.002801730910567093
$ws.Range("J14").Value = 0.04356334635268411
$ws.Range("L14").Value = 0.7303614371082858
$ws.Range("M14").Value = 1.706786105403566
$ws.Range("N14").Value = 6.580459192939884
$ws.Range("B15").Value = 9.325627705542502
$ws.Range("C15").Value = 1.131113059266056
$ws.Range("D15").Value = 0.2837111948989559
$ws.Range("E15").Value = 0.03835514933519057
$ws.Range("G15").Value = 0.002802729691346144
$ws.Range("J15").Value = 0.04310193070914892
$ws.Range("L15").Value = 0.7298610939391068
$ws.Range("M15").Value = 1.703784757903847
$ws.Range("N15").Value = 6.569890917739201
$ws.Range("B16").Value = 9.210900714776471
$ws.Range("C16").Value = 1.102504551295681
$ws.Range("D16").Value = 0.2838773254412743
$ws.Range("E16").Value = 0.03803549295332331
$ws.Range("G16").Value = 0.002808536190636383
$ws.Range("J16").Value = 0.0404632052519176
$ws.Range("L16").Value = 0.7270976011431713
$ws.Range("M16").Value = 1.686871303023182
$ws.Range("N16").Value = 6.509505241582531
$ws.Range("B17").Value = 9.141967681843425
$ws.Range("C17").Value = 1.085155580707465
$ws.Range("D17").Value = 0.2840065381418952
$ws.Range("E17").Value = 0.03783827085538327
$ws.Range("G17").Value = 0.002812172495196991
$ws.Range("J17").Value = 0.03884914317023913
$ws.Range("L17").Value = 0.7254946487154257
$ws.Range("M17").Value = 1.676749558425001
$ws.Range("N17").Value = 6.472615840411038
$ws.Range("B18").Value = 9.102848533973088
$ws.Range("C18").Value = 1.075250201410654
$ws.Range("D18").Value = 0.2840908846751518
$ws.Range("E18").Value = 0.03772441311077657
$ws.Range("G18").Value = 0.002814291366888384
$ws.Range("J18").Value = 0.03792244456408866
$ws.Range("L18").Value = 0.7246065229647058
$ws.Range("M18").Value = 1.671020737959367
$ws.Range("N18").Value = 6.451454088378739
$ws.Range("B19").Value = 9.089694098109476
$ws.Range("C19").Value = 1.071908928981259
$ws.Range("D19").Value = 0.2841211646495196
$ws.Range("E19").Value = 0.03768578980562509
$ws.Range("G19").Value = 0.002815013489726334
$ws.Range("J19").Value = 0.037608960675378
$ws.Range("L19").Value = 0.7243116240176306
$ws.Range("M19").Value = 1.669096978394364
$ws.Range("N19").Value = 6.444298616096319
$ws.Range("B20").Value = 9.149250874065388
$ws.Range("C20").Value = 1.086994806641769
$ws.Range("D20").Value = 0.2839917453911198
$ws.Range("E20").Value = 0.03785930885684863
$ws.Range("G20").Value = 0.002811782574360697
$ws.Range("J20").Value = 0.03902078885201377
$ws.Range("L20").Value = 0.7256617803274708
$ws.Range("M20").Value = 1.677817406585319
$ws.Range("N20").Value = 6.476536956137977
$ws.Range("B21").Value = 9.355698115596965
$ws.Range("C21").Value = 1.138562088876824
$ws.Range("D21").Value = 0.283676746514125
$ws.Range("E21").Value = 0.038437337643483
$ws.Range("G21").Value = 0.002801253431004896
$ws.Range("J21").Value = 0.04378470659714395
$ws.Range("L21").Value = 0.7306032042697694
$ws.Range("M21").Value = 1.708230409374949
$ws.Range("N21").Value = 6.585530132101553
$ws.Range("B22").Value = 9.495399526648725
$ws.Range("C22").Value = 1.172932631749802
$ws.Range("D22").Value = 0.2835602704982421
$ws.Range("E22").Value = 0.03881153699003193
$ws.Range("G22").Value = 0.00279461368941877
$ws.Range("J22").Value = 0.04691441645825023
$ws.Range("L22").Value = 0.7341361313258119
$ws.Range("M22").Value = 1.728944268102708
$ws.Range("N22").Value = 6.657284021989938
$ws.Range("B23").Value = 9.420398134226161
$ws.Range("C23").Value = 1.154526500056875
$ws.Range("D23").Value = 0.2836142507751021
$ws.Range("E23").Value = 0.03861213856107604
$ws.Range("G23").Value = 0.002798135430259281
$ws.Range("J23").Value = 0.04524248069606074
$ws.Range("L23").Value = 0.7322226797881086
$ws.Range("M23").Value = 1.717811814046797
$ws.Range("N23").Value = 6.618938944186937
$ws.Range("B24").Value = 9.145956550742994
$ws.Range("C24").Value = 1.086163078550896
$ws.Range("D24").Value = 0.2839984018584474
$ws.Range("E24").Value = 0.03784979904883912
$ws.Range("G24").Value = 0.002811958769362612
$ws.Range("J24").Value = 0.03894318396945806
$ws.Range("L24").Value = 0.7255861160111294
$ws.Range("M24").Value = 1.677334351256107
$ws.Range("N24").Value = 6.474764075252722
$ws.Range("B25").Value = 8.870127399907233
$ws.Range("C25").Value = 1.01521499275708
$ws.Range("D25").Value = 0.2847960683164317
$ws.Range("E25").Value = 0.03701130380776885
$ws.Range("G25").Value = 0.002827916946080832
$ws.Range("J25").Value = 0.03221131604166061
$ws.Range("L25").Value = 0.7197205487892546
$ws.Range("M25").Value = 1.637220624381726
$ws.Range("N25").Value = 6.321378459357845
